# Fruta / hortaliza, semanal
# Insert two new weekly price records for "Terminal Hortofrutícola Agro
# Chillán - Arándano (blue)" right after the first data row (row 2), pushing
# the rest of the historical rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at rows 3-4 (existing rows 3..22 shift to 5..24).
$ws.Rows("3:4").Insert()

# --- New row 3: Primera quality, 2023-02-06 ---
$ws.Cells.Item(3, 1).Value2 = 7
$ws.Cells.Item(3, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(3, 3).Value2 = "Ñuble"
$ws.Cells.Item(3, 4).Value2 = 44963
$ws.Cells.Item(3, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(3, 5).Value2 = 16
$ws.Cells.Item(3, 6).Value2 = "Fruta"
$ws.Cells.Item(3, 7).Value2 = 100101
$ws.Cells.Item(3, 8).Value2 = "Berries"
$ws.Cells.Item(3, 9).Value2 = 100101001
$ws.Cells.Item(3, 10).Value2 = "Arándano (blue)"
$ws.Cells.Item(3, 11).Value2 = "Sin especificar"
$ws.Cells.Item(3, 12).Value2 = "Primera"
$ws.Cells.Item(3, 13).Value2 = 50
$ws.Cells.Item(3, 14).Value2 = 3000
$ws.Cells.Item(3, 15).Value2 = 3000
$ws.Cells.Item(3, 16).Value2 = 3000
$ws.Cells.Item(3, 17).Value2 = "$/bandeja 2 kilos"
$ws.Cells.Item(3, 18).Value2 = "Provincia de Diguillín"
$ws.Cells.Item(3, 19).Value2 = 1500
$ws.Cells.Item(3, 20).Value2 = 2

# --- New row 4: Segunda quality, 2023-02-06 ---
$ws.Cells.Item(4, 1).Value2 = 7
$ws.Cells.Item(4, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(4, 3).Value2 = "Ñuble"
$ws.Cells.Item(4, 4).Value2 = 44963
$ws.Cells.Item(4, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(4, 5).Value2 = 16
$ws.Cells.Item(4, 6).Value2 = "Fruta"
$ws.Cells.Item(4, 7).Value2 = 100101
$ws.Cells.Item(4, 8).Value2 = "Berries"
$ws.Cells.Item(4, 9).Value2 = 100101001
$ws.Cells.Item(4, 10).Value2 = "Arándano (blue)"
$ws.Cells.Item(4, 11).Value2 = "Sin especificar"
$ws.Cells.Item(4, 12).Value2 = "Segunda"
$ws.Cells.Item(4, 13).Value2 = 50
$ws.Cells.Item(4, 14).Value2 = 2500
$ws.Cells.Item(4, 15).Value2 = 2500
$ws.Cells.Item(4, 16).Value2 = 2500
$ws.Cells.Item(4, 17).Value2 = "$/bandeja 2 kilos"
$ws.Cells.Item(4, 18).Value2 = "Provincia de Diguillín"
$ws.Cells.Item(4, 19).Value2 = 1250
$ws.Cells.Item(4, 20).Value2 = 2
